# Generate Report for Handoff
# The ca725fb2-... file has finished its round-trip (it was removed from the
# localization status report) and the remaining 911c866f-... entry moved
# from "Handed back" to "Ready for handoff" with refreshed timestamps.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview": drop row 3 (ca725fb2 file), refresh status/datetime
# for the remaining row (911c866f file).
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Rows.Item(3).Delete()

$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"
$wsOverview.Range("D2").Value = "2016-03-22 00:45:24"

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/77b5205bca2c75ea8a5c32f60d63e14eca5b8342/e2e/911c866f-3374-4099-b8f3-5a6de249a572.md", "", "", "911c866f-3374-4099-b8f3-5a6de249a572.md")

# ---------------------------------------------------------------------
# Sheet "zh-cn": drop row 3 (ca725fb2 file), refresh status/datetime for
# the remaining row (911c866f file).
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Rows.Item(3).Delete()

$wsZh.Range("C2").Value = "Ready for handoff"
$wsZh.Range("E2").Value = "2016-03-22 00:45:16"

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/77b5205bca2c75ea8a5c32f60d63e14eca5b8342/e2e/911c866f-3374-4099-b8f3-5a6de249a572.md", "", "", "911c866f-3374-4099-b8f3-5a6de249a572.md")
$wsZh.Hyperlinks.Add($wsZh.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/59f78849ba7db28d5acba908e5ff47917aa6931d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/911c866f-3374-4099-b8f3-5a6de249a572.cf28377ec91c04a903c845dbed1b68185b3cfd30.zh-cn.xlf", "", "", "911c866f-3374-4099-b8f3-5a6de249a572.cf28377ec91c04a903c845dbed1b68185b3cfd30.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/5aeb3cf89839cb0907649fb690b334021c47be02/e2e/911c866f-3374-4099-b8f3-5a6de249a572.md", "", "", "911c866f-3374-4099-b8f3-5a6de249a572.md")
$wsZh.Hyperlinks.Add($wsZh.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/bf3d4d2d64af1facb04ff54eda2571411354eb34/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/911c866f-3374-4099-b8f3-5a6de249a572.cf28377ec91c04a903c845dbed1b68185b3cfd30.zh-cn.xlf", "", "", "911c866f-3374-4099-b8f3-5a6de249a572.cf28377ec91c04a903c845dbed1b68185b3cfd30.zh-cn.xlf")

# ---------------------------------------------------------------------
# Sheet "de-de": drop row 3 (ca725fb2 file), refresh status/datetime for
# the remaining row (911c866f file).
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Rows.Item(3).Delete()

$wsDe.Range("C2").Value = "Ready for handoff"
$wsDe.Range("E2").Value = "2016-03-22 00:45:24"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/77b5205bca2c75ea8a5c32f60d63e14eca5b8342/e2e/911c866f-3374-4099-b8f3-5a6de249a572.md", "", "", "911c866f-3374-4099-b8f3-5a6de249a572.md")
$wsDe.Hyperlinks.Add($wsDe.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/653710fb5e4a66f72260fc8cd0f307de45122294/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/911c866f-3374-4099-b8f3-5a6de249a572.cf28377ec91c04a903c845dbed1b68185b3cfd30.de-de.xlf", "", "", "911c866f-3374-4099-b8f3-5a6de249a572.cf28377ec91c04a903c845dbed1b68185b3cfd30.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/53dbb55cbceb370f274944effc1a6a32db6669fa/e2e/911c866f-3374-4099-b8f3-5a6de249a572.md", "", "", "911c866f-3374-4099-b8f3-5a6de249a572.md")
$wsDe.Hyperlinks.Add($wsDe.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/84747f0629a0b45684b33a1ab8c10425a9d909b6/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/911c866f-3374-4099-b8f3-5a6de249a572.cf28377ec91c04a903c845dbed1b68185b3cfd30.de-de.xlf", "", "", "911c866f-3374-4099-b8f3-5a6de249a572.cf28377ec91c04a903c845dbed1b68185b3cfd30.de-de.xlf")
